$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (becomes sheet #2)
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$newws = $wb.Worksheets.Add($null, $wsTotal)
$newws.Name = "2022-Q3"

# Source sheet to copy cell formatting (bold/centered/bordered header style)
# from - fetch it *by name*, and *after* the Add() above, since sheet
# references resolve by live position/name rather than a frozen pointer.
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row format (B1:H1) and the index-column format (A-col)
# from the existing "2022-Q2" sheet onto the new sheet so the new sheet's
# style indices line up with the workbook's existing bold/center/border style.
$wsQ2.Range("B1:H1").Copy()
$newws.Range("B1:H1").PasteSpecial(-4122)

$wsQ2.Range("A2").Copy()
$newws.Range("A2:A20").PasteSpecial(-4122)

# Header row text
$newws.Cells.Item(1,2).Value = "基金代码"
$newws.Cells.Item(1,3).Value = "基金名称"
$newws.Cells.Item(1,4).Value = "基金规模"
$newws.Cells.Item(1,5).Value = "股票总仓位"
$newws.Cells.Item(1,6).Value = "仓位占比"
$newws.Cells.Item(1,7).Value = "持有市值(亿元)"
$newws.Cells.Item(1,8).Value = "仓位排名"

# Columns B-G (fund code / name / size / position% / ratio / market value)
# are stored as *text* in the source workbook (values like "32.70" or
# "002560" must keep their formatting / leading zeros), so force a text
# number format before writing the values.
$newws.Range("B2:G20").NumberFormat = "@"

$newws.Cells.Item(2,1).Value = 0
$newws.Cells.Item(2,2).Value = "002560"
$newws.Cells.Item(2,3).Value = "诺安和鑫灵活配置混合"
$newws.Cells.Item(2,4).Value = "32.70"
$newws.Cells.Item(2,5).Value = "79.56"
$newws.Cells.Item(2,6).Value = "7.60"
$newws.Cells.Item(2,7).Value = "2.4852"
$newws.Cells.Item(2,8).Value = 3
$newws.Cells.Item(3,1).Value = 1
$newws.Cells.Item(3,2).Value = "001822"
$newws.Cells.Item(3,3).Value = "华商智能生活灵活配置混合A"
$newws.Cells.Item(3,4).Value = "33.45"
$newws.Cells.Item(3,5).Value = "87.34"
$newws.Cells.Item(3,6).Value = "6.80"
$newws.Cells.Item(3,7).Value = "2.2746"
$newws.Cells.Item(3,8).Value = 3
$newws.Cells.Item(4,1).Value = 2
$newws.Cells.Item(4,2).Value = "001933"
$newws.Cells.Item(4,3).Value = "华商新兴活力灵活配置混合"
$newws.Cells.Item(4,4).Value = "18.34"
$newws.Cells.Item(4,5).Value = "87.96"
$newws.Cells.Item(4,6).Value = "6.87"
$newws.Cells.Item(4,7).Value = "1.2600"
$newws.Cells.Item(4,8).Value = 3
$newws.Cells.Item(5,1).Value = 3
$newws.Cells.Item(5,2).Value = "010550"
$newws.Cells.Item(5,3).Value = "华商双擎领航混合"
$newws.Cells.Item(5,4).Value = "12.53"
$newws.Cells.Item(5,5).Value = "90.32"
$newws.Cells.Item(5,6).Value = "6.78"
$newws.Cells.Item(5,7).Value = "0.8495"
$newws.Cells.Item(5,8).Value = 4
$newws.Cells.Item(6,1).Value = 4
$newws.Cells.Item(6,2).Value = "015385"
$newws.Cells.Item(6,3).Value = "华商智能生活灵活配置混合C"
$newws.Cells.Item(6,4).Value = "11.97"
$newws.Cells.Item(6,5).Value = "87.34"
$newws.Cells.Item(6,6).Value = "6.80"
$newws.Cells.Item(6,7).Value = "0.8140"
$newws.Cells.Item(6,8).Value = 3
$newws.Cells.Item(7,1).Value = 5
$newws.Cells.Item(7,2).Value = "013886"
$newws.Cells.Item(7,3).Value = "华商新能源汽车混合A"
$newws.Cells.Item(7,4).Value = "9.72"
$newws.Cells.Item(7,5).Value = "89.01"
$newws.Cells.Item(7,6).Value = "6.99"
$newws.Cells.Item(7,7).Value = "0.6794"
$newws.Cells.Item(7,8).Value = 4
$newws.Cells.Item(8,1).Value = 6
$newws.Cells.Item(8,2).Value = "001411"
$newws.Cells.Item(8,3).Value = "诺安创新驱动灵活配置混合A"
$newws.Cells.Item(8,4).Value = "3.98"
$newws.Cells.Item(8,5).Value = "80.56"
$newws.Cells.Item(8,6).Value = "7.37"
$newws.Cells.Item(8,7).Value = "0.2933"
$newws.Cells.Item(8,8).Value = 5
$newws.Cells.Item(9,1).Value = 7
$newws.Cells.Item(9,2).Value = "013887"
$newws.Cells.Item(9,3).Value = "华商新能源汽车混合C"
$newws.Cells.Item(9,4).Value = "4.06"
$newws.Cells.Item(9,5).Value = "89.01"
$newws.Cells.Item(9,6).Value = "6.99"
$newws.Cells.Item(9,7).Value = "0.2838"
$newws.Cells.Item(9,8).Value = 4
$newws.Cells.Item(10,1).Value = 8
$newws.Cells.Item(10,2).Value = "001239"
$newws.Cells.Item(10,3).Value = "长盛国企改革主题灵活配置混合"
$newws.Cells.Item(10,4).Value = "4.46"
$newws.Cells.Item(10,5).Value = "90.97"
$newws.Cells.Item(10,6).Value = "4.93"
$newws.Cells.Item(10,7).Value = "0.2199"
$newws.Cells.Item(10,8).Value = 8
$newws.Cells.Item(11,1).Value = 9
$newws.Cells.Item(11,2).Value = "014350"
$newws.Cells.Item(11,3).Value = "华商卓越成长一年持有混合A"
$newws.Cells.Item(11,4).Value = "3.14"
$newws.Cells.Item(11,5).Value = "86.88"
$newws.Cells.Item(11,6).Value = "6.84"
$newws.Cells.Item(11,7).Value = "0.2148"
$newws.Cells.Item(11,8).Value = 3
$newws.Cells.Item(12,1).Value = 10
$newws.Cells.Item(12,2).Value = "010852"
$newws.Cells.Item(12,3).Value = "中欧内需成长混合A"
$newws.Cells.Item(12,4).Value = "3.70"
$newws.Cells.Item(12,5).Value = "90.11"
$newws.Cells.Item(12,6).Value = "5.15"
$newws.Cells.Item(12,7).Value = "0.1906"
$newws.Cells.Item(12,8).Value = 8
$newws.Cells.Item(13,1).Value = 11
$newws.Cells.Item(13,2).Value = "002051"
$newws.Cells.Item(13,3).Value = "诺安创新驱动灵活配置混合C"
$newws.Cells.Item(13,4).Value = "1.80"
$newws.Cells.Item(13,5).Value = "80.56"
$newws.Cells.Item(13,6).Value = "7.37"
$newws.Cells.Item(13,7).Value = "0.1327"
$newws.Cells.Item(13,8).Value = 5
$newws.Cells.Item(14,1).Value = 12
$newws.Cells.Item(14,2).Value = "970043"
$newws.Cells.Item(14,3).Value = "东吴裕盈一年持有期灵活配置混合A"
$newws.Cells.Item(14,4).Value = "0.96"
$newws.Cells.Item(14,5).Value = "52.43"
$newws.Cells.Item(14,6).Value = "5.79"
$newws.Cells.Item(14,7).Value = "0.0556"
$newws.Cells.Item(14,8).Value = 2
$newws.Cells.Item(15,1).Value = 13
$newws.Cells.Item(15,2).Value = "010853"
$newws.Cells.Item(15,3).Value = "中欧内需成长混合C"
$newws.Cells.Item(15,4).Value = "0.54"
$newws.Cells.Item(15,5).Value = "90.11"
$newws.Cells.Item(15,6).Value = "5.15"
$newws.Cells.Item(15,7).Value = "0.0278"
$newws.Cells.Item(15,8).Value = 8
$newws.Cells.Item(16,1).Value = 14
$newws.Cells.Item(16,2).Value = "970045"
$newws.Cells.Item(16,3).Value = "东吴裕盈一年持有期灵活配置混合C"
$newws.Cells.Item(16,4).Value = "0.44"
$newws.Cells.Item(16,5).Value = "52.43"
$newws.Cells.Item(16,6).Value = "5.79"
$newws.Cells.Item(16,7).Value = "0.0255"
$newws.Cells.Item(16,8).Value = 2
$newws.Cells.Item(17,1).Value = 15
$newws.Cells.Item(17,2).Value = "970044"
$newws.Cells.Item(17,3).Value = "东吴裕盈一年持有期灵活配置混合B"
$newws.Cells.Item(17,4).Value = "0.27"
$newws.Cells.Item(17,5).Value = "52.43"
$newws.Cells.Item(17,6).Value = "5.79"
$newws.Cells.Item(17,7).Value = "0.0156"
$newws.Cells.Item(17,8).Value = 2
$newws.Cells.Item(18,1).Value = 16
$newws.Cells.Item(18,2).Value = "168701"
$newws.Cells.Item(18,3).Value = "合煦智远国证香蜜湖金融科技指数（LOF）A"
$newws.Cells.Item(18,4).Value = "0.47"
$newws.Cells.Item(18,5).Value = "92.07"
$newws.Cells.Item(18,6).Value = "2.13"
$newws.Cells.Item(18,7).Value = "0.0100"
$newws.Cells.Item(18,8).Value = 8
$newws.Cells.Item(19,1).Value = 17
$newws.Cells.Item(19,2).Value = "014351"
$newws.Cells.Item(19,3).Value = "华商卓越成长一年持有混合C"
$newws.Cells.Item(19,4).Value = "0.10"
$newws.Cells.Item(19,5).Value = "86.88"
$newws.Cells.Item(19,6).Value = "6.84"
$newws.Cells.Item(19,7).Value = "0.0068"
$newws.Cells.Item(19,8).Value = 3
$newws.Cells.Item(20,1).Value = 18
$newws.Cells.Item(20,2).Value = "168702"
$newws.Cells.Item(20,3).Value = "合煦智远国证香蜜湖金融科技指数（LOF）C"
$newws.Cells.Item(20,4).Value = "0.13"
$newws.Cells.Item(20,5).Value = "92.07"
$newws.Cells.Item(20,6).Value = "2.13"
$newws.Cells.Item(20,7).Value = "0.0028"
$newws.Cells.Item(20,8).Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 and
#    shift the existing history down by one row, bumping the index
#    column (A) for every shifted row.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# The freshly inserted row 2 inherits row 1's (header) bold/center format;
# the data rows in this sheet are unstyled, so strip that back off.
$wsTotal.Range("B2:D2").ClearFormats()

# A2 should carry the same bold/center/bordered style as the rest of the
# index column (A3 now holds that formatting, shifted down from old A2).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 19
$wsTotal.Cells.Item(2,4).Value = 9.84

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(6,1).Value = 4
$wsTotal.Cells.Item(7,1).Value = 5
$wsTotal.Cells.Item(8,1).Value = 6
$wsTotal.Cells.Item(9,1).Value = 7
